$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial number. All rows from 2 to 257
# had the value 45177 (2023-09-08), which needs to be updated to 45178
# (2023-09-09).
$ws.Range("C2:C257").Value = 45178
